{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"2025-11-12 Wednesday\", \"2025-11-13 Thursday\"],\n  [\"466\u00d77=\", \"321\u00d72=\"],\n  [\"468\u00d76=\", \"511\u00d75=\"],\n  [\"475\u00d79=\", \"340\u00d75=\"],\n  [\"104\u00d76=\", \"938\u00d73=\"],\n  [\"882\u00d72=\", \"944\u00d78=\"],\n  [\"804\u00d72=\", \"782\u00d76=\"],\n  [\"668\u00d79=\", \"899\u00d78=\"],\n  [\"166\u00d72=\", \"966\u00d77=\"],\n  [\"455\u00d78=\", \"907\u00d75=\"],\n  [\"415\u00d76=\", \"655\u00d75=\"],\n  [\"135\u00d77=\", \"537\u00d77=\"],\n  [\"283\u00d74=\", \"466\u00d74=\"],\n  [\"143\u00d76=\", \"426\u00d75=\"],\n  [\"412\u00d76=\", \"401\u00d79=\"],\n  [\"925\u00d75=\", \"201\u00d79=\"],\n  [\"999\u00d78=\", \"559\u00d77=\"],\n  [\"783\u00d76=\", \"940\u00d78=\"],\n  [\"219\u00d73=\", \"322\u00d75=\"],\n  [\"980\u00d74=\", \"623\u00d79=\"],\n  [\"695\u00d77=\", \"618\u00d77=\"],\n  [\"734\u00d77=\", \"296\u00d79=\"],\n  [\"425\u00d76=\", \"885\u00d77=\"],\n  [\"911\u00d74=\", \"114\u00d77=\"],\n  [\"977\u00d74=\", \"965\u00d78=\"],\n  [\"790\u00d75=\", \"277\u00d79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"2025-11-12 Wednesday\"\n$find.Replacement.Text = \"2025-11-13 Thursday\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"466\u00d77=\"\n$find.Replacement.Text = \"321\u00d72=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"468\u00d76=\"\n$find.Replacement.Text = \"511\u00d75=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"475\u00d79=\"\n$find.Replacement.Text = \"340\u00d75=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"104\u00d76=\"\n$find.Replacement.Text = \"938\u00d73=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"882\u00d72=\"\n$find.Replacement.Text = \"944\u00d78=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"804\u00d72=\"\n$find.Replacement.Text = \"782\u00d76=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"668\u00d79=\"\n$find.Replacement.Text = \"899\u00d78=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"166\u00d72=\"\n$find.Replacement.Text = \"966\u00d77=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"455\u00d78=\"\n$find.Replacement.Text = \"907\u00d75=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"415\u00d76=\"\n$find.Replacement.Text = \"655\u00d75=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"135\u00d77=\"\n$find.Replacement.Text = \"537\u00d77=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"283\u00d74=\"\n$find.Replacement.Text = \"466\u00d74=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"143\u00d76=\"\n$find.Replacement.Text = \"426\u00d75=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"412\u00d76=\"\n$find.Replacement.Text = \"401\u00d79=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"925\u00d75=\"\n$find.Replacement.Text = \"201\u00d79=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"999\u00d78=\"\n$find.Replacement.Text = \"559\u00d77=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"783\u00d76=\"\n$find.Replacement.Text = \"940\u00d78=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"219\u00d73=\"\n$find.Replacement.Text = \"322\u00d75=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"980\u00d74=\"\n$find.Replacement.Text = \"623\u00d79=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"695\u00d77=\"\n$find.Replacement.Text = \"618\u00d77=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"734\u00d77=\"\n$find.Replacement.Text = \"296\u00d79=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"425\u00d76=\"\n$find.Replacement.Text = \"885\u00d77=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"911\u00d74=\"\n$find.Replacement.Text = \"114\u00d77=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"977\u00d74=\"\n$find.Replacement.Text = \"965\u00d78=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"790\u00d75=\"\n$find.Replacement.Text = \"277\u00d79=\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n"}
